# "Add files via upload" — the uploaded "Price Target.xlsx" now reflects a
# fresh quote pull: the helper column B labels were renamed from
# "target*" to "target*Price" (targetHigh -> targetHighPrice, etc.) and
# that (hidden) helper column was widened. The live xlquotePriceTarget()
# array-formula results in column C are left exactly as-is (formula,
# array-ref and dynamic-array cell metadata untouched) since those values
# come from the external XLL quote add-in, not something this script
# recomputes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the lookup-metric helper cells (column B) that back each
# _xll.xlquotePriceTarget(Symbol, B#) array formula in column C.
$ws.Range("B4").Value = "targetHighPrice"
$ws.Range("B5").Value = "targetLowPrice"
$ws.Range("B6").Value = "targetMeanPrice"
$ws.Range("B7").Value = "targetMedianPrice"

# The helper column stays hidden but is now given an explicit width.
$ws.Range("B1").EntireColumn.ColumnWidth = 14.59
